$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: replace ASUSTeK supplier invoice with OpenERP invoice, shift dates back 2 years (1 day earlier ref)
$ws.Range("A3").Value = "PI12/0101"
$ws.Range("B3").Value = "OpenERP"
$ws.Range("D3").Value = 39824
$ws.Range("H3").Value = "Open Invoice from supplier OpenERP"

# Row 4: update reference number and date, wording change supplier -> customer
$ws.Range("A4").Value = "PI12/0250"
$ws.Range("B4").Value = "032/0029"
$ws.Range("D4").Value = 39844
$ws.Range("H4").Value = "Open Invoice from customer with reference '032/0029'"

# Update active cell selection on the sheet view to K7
$ws.Range("K7").Select()
